$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value2 = 1580.4117
$ws.Range("I28").Value2 = 2557.2222
$ws.Range("J28").Value2 = 481.5
$ws.Range("K28").Value2 = 2557.2222
$ws.Range("L28").Value2 = 481.5
$ws.Range("M28").Value2 = -2072.2222
$ws.Range("N28").Value2 = -1451.5
# Row 103
$ws.Range("H103").Value2 = 774.48
$ws.Range("I103").Value2 = 857.2857
$ws.Range("J103").Value2 = 669.0909
$ws.Range("K103").Value2 = 2571.8571
$ws.Range("L103").Value2 = 2007.2727
$ws.Range("M103").Value2 = -1985.8571
$ws.Range("N103").Value2 = -3179.2727
# Row 129
$ws.Range("H129").Value2 = 43210756
$ws.Range("J129").Value2 = 2179664
$ws.Range("L129").Value2 = 6538992
$ws.Range("N129").Value2 = -6548992

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value2 = 15117.56
$ws.Range("I32").Value2 = 13181.681
$ws.Range("J32").Value2 = 45446.332
$ws.Range("K32").Value2 = 13181.681
$ws.Range("L32").Value2 = 45446.332
$ws.Range("M32").Value2 = -12894.681
$ws.Range("N32").Value2 = -46020.332
# Row 86
$ws.Range("H86").Value2 = 37995
$ws.Range("J86").Value2 = 37995
$ws.Range("L86").Value2 = 37995
$ws.Range("N86").Value2 = -40367
# Row 89
$ws.Range("H89").Value2 = 37995
$ws.Range("J89").Value2 = 37995
$ws.Range("L89").Value2 = 113985
$ws.Range("N89").Value2 = -125841
# Row 132
$ws.Range("H132").Value2 = 2448
$ws.Range("I132").Value2 = 2064.3333
$ws.Range("J132").Value2 = 2754.9333
$ws.Range("K132").Value2 = 6192.999899999999
$ws.Range("L132").Value2 = 8264.7999
$ws.Range("M132").Value2 = -3662.999899999999
$ws.Range("N132").Value2 = -13324.7999

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 88
$ws.Range("H88").Value2 = 0
$ws.Range("J88").Value2 = 0
$ws.Range("L88").Value2 = 0
$ws.Range("N88").ClearContents() | Out-Null
# Row 91
$ws.Range("H91").Value2 = 0
$ws.Range("J91").Value2 = 0
$ws.Range("L91").Value2 = 0
$ws.Range("N91").ClearContents() | Out-Null
# Row 134
$ws.Range("H134").Value2 = 34431.547
$ws.Range("I134").Value2 = 1895.12
$ws.Range("J134").Value2 = 170000
$ws.Range("K134").Value2 = 5685.36
$ws.Range("L134").Value2 = 510000
$ws.Range("M134").Value2 = -3150.36
$ws.Range("N134").Value2 = -515070

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 88
$ws.Range("H88").Value2 = 35000
$ws.Range("J88").Value2 = 35000
$ws.Range("L88").Value2 = 35000
$ws.Range("N88").Value2 = -35812
# Row 91
$ws.Range("H91").Value2 = 35000
$ws.Range("J91").Value2 = 35000
$ws.Range("L91").Value2 = 35000
$ws.Range("N91").Value2 = -37808
# Row 99
$ws.Range("H99").Value2 = 2277.8262
$ws.Range("I99").Value2 = 1869.5
$ws.Range("K99").Value2 = 1869.5
$ws.Range("M99").Value2 = -371.5
# Row 122
$ws.Range("H122").Value2 = 909970.4
$ws.Range("I122").Value2 = 1429469.1
$ws.Range("J122").Value2 = 847.5
$ws.Range("K122").Value2 = 4288407.300000001
$ws.Range("L122").Value2 = 2542.5
$ws.Range("M122").Value2 = -4285957.300000001
$ws.Range("N122").Value2 = -7442.5
# Row 126
$ws.Range("H126").Value2 = 2277.8262
$ws.Range("I126").Value2 = 1869.5
$ws.Range("K126").Value2 = 5608.5
$ws.Range("M126").Value2 = -3138.5
# Row 134
$ws.Range("H134").Value2 = 1871.3125
$ws.Range("I134").Value2 = 1256.1471
$ws.Range("J134").Value2 = 3365.2856
$ws.Range("K134").Value2 = 3768.4413
$ws.Range("L134").Value2 = 10095.8568
$ws.Range("M134").Value2 = -1233.4413
$ws.Range("N134").Value2 = -15165.8568

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value2 = 1335.5781
$ws.Range("I68").Value2 = 710.125
$ws.Range("J68").Value2 = 1961.0312
$ws.Range("K68").Value2 = 2130.375
$ws.Range("L68").Value2 = 5883.0936
$ws.Range("M68").Value2 = -1319.375
$ws.Range("N68").Value2 = -7505.0936
# Row 71
$ws.Range("H71").Value2 = 1335.5781
$ws.Range("I71").Value2 = 710.125
$ws.Range("J71").Value2 = 1961.0312
$ws.Range("K71").Value2 = 6391.125
$ws.Range("L71").Value2 = 17649.2808
$ws.Range("M71").Value2 = -2335.125
$ws.Range("N71").Value2 = -25761.2808
# Row 75
$ws.Range("H75").Value2 = 533
$ws.Range("I75").Value2 = 549.5
$ws.Range("J75").Value2 = 500
$ws.Range("K75").Value2 = 1648.5
$ws.Range("L75").Value2 = 1500
$ws.Range("M75").Value2 = -650.5
$ws.Range("N75").Value2 = -3496
# Row 78
$ws.Range("H78").Value2 = 533
$ws.Range("I78").Value2 = 549.5
$ws.Range("J78").Value2 = 500
$ws.Range("K78").Value2 = 4945.5
$ws.Range("L78").Value2 = 4500
$ws.Range("M78").Value2 = 46.5
$ws.Range("N78").Value2 = -14484
# Row 113
$ws.Range("H113").Value2 = 405.05
$ws.Range("I113").Value2 = 371.08694
$ws.Range("J113").Value2 = 415.1948
$ws.Range("K113").Value2 = 1113.26082
$ws.Range("L113").Value2 = 1245.5844
$ws.Range("M113").Value2 = 1056.73918
$ws.Range("N113").Value2 = -5585.5844
# Row 114
$ws.Range("H114").Value2 = 1197.591
$ws.Range("I114").Value2 = 681.7
$ws.Range("J114").Value2 = 1627.5
$ws.Range("K114").Value2 = 2045.1
$ws.Range("L114").Value2 = 4882.5
$ws.Range("M114").Value2 = 1208.9
$ws.Range("N114").Value2 = -11390.5
# Row 131
$ws.Range("H131").Value2 = 19790.328
$ws.Range("J131").Value2 = 1902.7441
$ws.Range("L131").Value2 = 5708.2323
$ws.Range("N131").Value2 = -15788.2323
# Row 132
$ws.Range("H132").Value2 = 840
$ws.Range("J132").Value2 = 900
$ws.Range("L132").Value2 = 8100
$ws.Range("N132").Value2 = -13160

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value2 = 1210.8422
$ws.Range("I113").Value2 = 994
$ws.Range("J113").Value2 = 1451.7778
$ws.Range("K113").Value2 = 994
$ws.Range("L113").Value2 = 1451.7778
$ws.Range("M113").Value2 = 1176
$ws.Range("N113").Value2 = -5791.7778

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value2 = 11851232
$ws.Range("I40").Value2 = 13746973
$ws.Range("J40").Value2 = 2849.75
$ws.Range("K40").Value2 = 13746973
$ws.Range("M40").Value2 = -13746837
$ws.Range("N40").Value2 = -3121.75
# Row 80
$ws.Range("H80").Value2 = 19562.5
$ws.Range("J80").Value2 = 19562.5
$ws.Range("L80").Value2 = 19562.5
$ws.Range("N80").Value2 = -21808.5
# Row 83
$ws.Range("H83").Value2 = 19562.5
$ws.Range("J83").Value2 = 19562.5
$ws.Range("L83").Value2 = 58687.5
$ws.Range("N83").Value2 = -69919.5

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 87
$ws.Range("H87").Value2 = 0
$ws.Range("J87").Value2 = 0
$ws.Range("L87").Value2 = 0
$ws.Range("N87").ClearContents() | Out-Null
# Row 90
$ws.Range("H90").Value2 = 0
$ws.Range("J90").Value2 = 0
$ws.Range("L90").Value2 = 0
$ws.Range("N90").ClearContents() | Out-Null
# Row 107
$ws.Range("H107").Value2 = 426.54544
$ws.Range("I107").Value2 = 424.875
$ws.Range("J107").Value2 = 431
$ws.Range("K107").Value2 = 1274.625
$ws.Range("L107").Value2 = 1293
$ws.Range("M107").Value2 = 645.375
$ws.Range("N107").Value2 = -5133
# Row 132
$ws.Range("H132").Value2 = 1508.7222
$ws.Range("I132").Value2 = 779.2083
$ws.Range("J132").Value2 = 2967.75
$ws.Range("K132").Value2 = 2337.6249
$ws.Range("L132").Value2 = 8903.25
$ws.Range("M132").Value2 = 192.3751000000002
$ws.Range("N132").Value2 = -13963.25

